$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 18 entirely (delete the row, not just clear contents)
$ws.Rows("18:18").Delete()
